$wb = $excel.ActiveWorkbook

# The "想去人数" (number of interested attendees) counts increased for two
# events, and these values are duplicated across the "展览" and "全部类型"
# sheets, which mirror the same data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 350
    $ws.Range("F5").Value = 109
}
